$d = $word.ActiveDocument

# This revision appends three new paragraphs to the very end of the body
# (right after the existing last paragraph, "... tømmer databasen.", and
# before the section break):
#
#   <blank paragraph>
#   Karma
#   Vi har ikke været i stand til at få de eksisterende karma test til at
#   virke såvel som at skrive relevante nye test.

# 1) A blank paragraph right after the current last paragraph.
$tail = $d.Content
$tail.Collapse(0)
$tail.InsertParagraphAfter()

# 2) A paragraph containing "Karma".
$tail = $d.Content
$tail.Collapse(0)
$tail.InsertParagraphAfter()
$karmaPara = $d.Paragraphs.Last
$karmaStart = $d.Range($karmaPara.Range.Start, $karmaPara.Range.Start)
$karmaStart.InsertAfter("Karma")

# 3) A paragraph with the closing remark about the karma tests.
$tail = $d.Content
$tail.Collapse(0)
$tail.InsertParagraphAfter()
$notePara = $d.Paragraphs.Last
$noteStart = $d.Range($notePara.Range.Start, $notePara.Range.Start)
$noteStart.InsertAfter("Vi har ikke været i stand til at få de eksisterende karma test til at virke såvel som at skrive relevante nye test.")
